$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Gilberto-Mec. Tec. Res. Mat."
$ws.Range("D3").Value = "-"

$ws.Range("B4").Value = "Gilberto-Mec. Tec. Res. Mat."
$ws.Range("C4").Value = "[Elaine-Metalografia-2B, -, Elaine-Metalografia-2B, -]"
$ws.Range("D4").Value = "[Victor S.-Ajustagem-2B, Carlos-Tornearia-2B, Victor S.-Ajustagem-2B, Carlos-Tornearia-2B]"
$ws.Range("F4").Value = "Maria Celeste-Maq. Term. Fluxo"

$ws.Range("B6").Value = "Gilberto-Mec. Tec. Res. Mat."
$ws.Range("C6").Value = "[Elaine-Metalografia-2B, Emerson-Elet. Digi. Básica-2B, Elaine-Metalografia-2B, Emerson-Elet. Digi. Básica-2B]"
$ws.Range("D6").Value = "[Elcio D.-Des. Maq. CAD-T2-2B, Carlos-Tornearia-2B, Elcio D.-Des. Maq. CAD-T2-2B, Carlos-Tornearia-2B]"
$ws.Range("F6").Value = "Maria Celeste-Maq. Term. Fluxo"

$ws.Range("C7").Value = "[Victor S.-Ajustagem-2B, Emerson-Elet. Digi. Básica-2B, Victor S.-Ajustagem-2B, Emerson-Elet. Digi. Básica-2B]"
$ws.Range("D7").Value = "[Elcio D.-Des. Maq. CAD-T1-2B, Elcio D.-Des. Maq. CAD-T1-2B, Elcio D.-Des. Maq. CAD-T2-2B, Elcio D.-Des. Maq. CAD-T1-2B]"
